$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3186.8667
$ws.Range("I40").Value = 4037.625
$ws.Range("J40").Value = 2214.5715
$ws.Range("K40").Value = 4037.625
$ws.Range("L40").Value = 2214.5715
$ws.Range("M40").Value = -3862.625
$ws.Range("N40").Value = -2564.5715
$ws.Range("H42").Value = 188.16667
$ws.Range("I42").Value = 76.333336
$ws.Range("K42").Value = 229.000008
$ws.Range("M42").Value = 0.9999919999999918
$ws.Range("H76").Value = 3178083.8
$ws.Range("I76").Value = 3971623
$ws.Range("J76").Value = 3926.8572
$ws.Range("K76").Value = 3971623
$ws.Range("L76").Value = 3926.8572
$ws.Range("M76").Value = -3971308
$ws.Range("N76").Value = -4556.8572
$ws.Range("H79").Value = 3178083.8
$ws.Range("I79").Value = 3971623
$ws.Range("J79").Value = 3926.8572
$ws.Range("K79").Value = 3971623
$ws.Range("L79").Value = 3926.8572
$ws.Range("M79").Value = -3970531
$ws.Range("N79").Value = -6110.8572
$ws.Range("H125").Value = 8008453
$ws.Range("I125").Value = 577
$ws.Range("J125").Value = 14014360
$ws.Range("K125").Value = 5193
$ws.Range("L125").Value = 126129240
$ws.Range("M125").Value = -2733
$ws.Range("N125").Value = -126134160
$ws.Range("H129").Value = 1155.35
$ws.Range("J129").Value = 1192.7894
$ws.Range("L129").Value = 3578.3682
$ws.Range("N129").Value = -13578.3682
$ws.Range("H132").Value = 490342.75
$ws.Range("I132").Value = 607592.1
$ws.Range("J132").Value = 21345.2
$ws.Range("K132").Value = 1822776.3
$ws.Range("L132").Value = 64035.60000000001
$ws.Range("M132").Value = -1820246.3
$ws.Range("N132").Value = -69095.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5391.4443
$ws.Range("J122").Value = 11733.333
$ws.Range("L122").Value = 35199.999
$ws.Range("N122").Value = -40099.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 19183.666
$ws.Range("I86").Value = 3019
$ws.Range("K86").Value = 3019
$ws.Range("M86").Value = -1896
$ws.Range("H89").Value = 19183.666
$ws.Range("I89").Value = 3019
$ws.Range("K89").Value = 15095
$ws.Range("M89").Value = -9479
$ws.Range("H134").Value = 3774.0527
$ws.Range("I134").Value = 3105.1428
$ws.Range("J134").Value = 5647
$ws.Range("K134").Value = 9315.428400000001
$ws.Range("L134").Value = 16941
$ws.Range("M134").Value = -6780.428400000001
$ws.Range("N134").Value = -22011

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H62").Value = 16369.4375
$ws.Range("I62").Value = 26377.777
$ws.Range("J62").Value = 3501.5715
$ws.Range("K62").Value = 26377.777
$ws.Range("L62").Value = 3501.5715
$ws.Range("M62").Value = -25753.777
$ws.Range("N62").Value = -4749.5715
$ws.Range("H65").Value = 16369.4375
$ws.Range("I65").Value = 26377.777
$ws.Range("J65").Value = 3501.5715
$ws.Range("K65").Value = 131888.885
$ws.Range("L65").Value = 17507.8575
$ws.Range("M65").Value = -128768.885
$ws.Range("N65").Value = -23747.8575
$ws.Range("H99").Value = 4471785
$ws.Range("I99").Value = 6955021
$ws.Range("J99").Value = 1959.8
$ws.Range("K99").Value = 6955021
$ws.Range("L99").Value = 1959.8
$ws.Range("M99").Value = -6953523
$ws.Range("N99").Value = -4955.8
$ws.Range("H126").Value = 4471785
$ws.Range("I126").Value = 6955021
$ws.Range("J126").Value = 1959.8
$ws.Range("K126").Value = 20865063
$ws.Range("L126").Value = 5879.4
$ws.Range("M126").Value = -20862593
$ws.Range("N126").Value = -10819.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4167867.5
$ws.Range("I131").Value = 213.6
$ws.Range("J131").Value = 4445711
$ws.Range("K131").Value = 640.8
$ws.Range("L131").Value = 13337133
$ws.Range("M131").Value = 4399.2
$ws.Range("N131").Value = -13347213

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6675.3213
$ws.Range("I70").Value = 8061.8184
$ws.Range("J70").Value = 5778.1763
$ws.Range("K70").Value = 8061.8184
$ws.Range("L70").Value = 5778.1763
$ws.Range("M70").Value = -7791.8184
$ws.Range("N70").Value = -6318.1763
$ws.Range("H73").Value = 6675.3213
$ws.Range("I73").Value = 8061.8184
$ws.Range("J73").Value = 5778.1763
$ws.Range("K73").Value = 8061.8184
$ws.Range("L73").Value = 5778.1763
$ws.Range("M73").Value = -7125.8184
$ws.Range("N73").Value = -7650.1763
$ws.Range("H107").Value = 484.85715
$ws.Range("I107").Value = 478.2
$ws.Range("J107").Value = 501.5
$ws.Range("K107").Value = 478.2
$ws.Range("L107").Value = 501.5
$ws.Range("M107").Value = 1441.8
$ws.Range("N107").Value = -4341.5
$ws.Range("H122").Value = 696516.0600000001
$ws.Range("I122").Value = 1235790.1
$ws.Range("J122").Value = 3163.7144
$ws.Range("K122").Value = 3707370.3
$ws.Range("L122").Value = 9491.143199999999
$ws.Range("M122").Value = -3704920.3
$ws.Range("N122").Value = -14391.1432
$ws.Range("H123").Value = 10956.315
$ws.Range("J123").Value = 10956.315
$ws.Range("L123").Value = 10956.315
$ws.Range("N123").Value = -15856.315
$ws.Range("H136").Value = 12872.333
$ws.Range("J136").Value = 12872.333
$ws.Range("L136").Value = 38616.999
$ws.Range("N136").Value = -43716.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 20799.6
$ws.Range("I22").Value = 999.6667
$ws.Range("J22").Value = 50499.5
$ws.Range("K22").Value = 999.6667
$ws.Range("L22").Value = 50499.5
$ws.Range("M22").Value = -704.6667
$ws.Range("N22").Value = -51089.5
$ws.Range("H27").Value = 20799.6
$ws.Range("I27").Value = 999.6667
$ws.Range("J27").Value = 50499.5
$ws.Range("K27").Value = 999.6667
$ws.Range("L27").Value = 50499.5
$ws.Range("M27").Value = -892.6667
$ws.Range("N27").Value = -50713.5
$ws.Range("H93").Value = 1562.3334
$ws.Range("I93").Value = 1249.7142
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 1249.7142
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = -1.714199999999892
$ws.Range("N93").Value = -4496
$ws.Range("H122").Value = 3805.4285
$ws.Range("I122").Value = 3301
$ws.Range("J122").Value = 3924.1177
$ws.Range("K122").Value = 9903
$ws.Range("L122").Value = 11772.3531
$ws.Range("M122").Value = -7453
$ws.Range("N122").Value = -16672.3531

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 112350.11
$ws.Range("I122").Value = 125993.875
$ws.Range("J122").Value = 3200
$ws.Range("K122").Value = 377981.625
$ws.Range("L122").Value = 9600
$ws.Range("M122").Value = -375531.625
$ws.Range("N122").Value = -14500
$ws.Range("H137").Value = 61376.8
$ws.Range("J137").Value = 61376.8
$ws.Range("L137").Value = 61376.8
$ws.Range("N137").Value = -71576.8
